# "Yêu cầu của Thầy" — mark rows 28-31 and 44-45 as "Hoàn thành" (E = 100%)
# and (re)assign "Phân công" (F) to "Bùi, Kiều" for rows 28-31.
# Row 28 previously had a stray G-column note ("Bùi") that is replaced by the
# proper F-column assignment; rows 44-45 already had the right G-column note
# ("Kiều") so only the completion flag changes there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 28: clear the stray G28 note, mark done, assign F28.
$ws.Range("G28").Clear()
$ws.Range("E28").Value = 1
$ws.Range("F28").Value = "Bùi, Kiều"

# Rows 29-31: mark done, assign F column.
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = "Bùi, Kiều"

$ws.Range("E30").Value = 1
$ws.Range("F30").Value = "Bùi, Kiều"

$ws.Range("E31").Value = 1
$ws.Range("F31").Value = "Bùi, Kiều"

# Rows 44-45: just mark done, G column assignment already correct.
$ws.Range("E44").Value = 1
$ws.Range("E45").Value = 1

# Update the view state: scroll so row 13 is at the top and select G30.
$ws.Range("A13").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G30").Select()
